$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column, matching the style of the existing header row (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data values for the "Save" column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
